$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 147.75
$ws.Range("I6").Value = 147.75
$ws.Range("K6").Value = 443.25
$ws.Range("M6").Value = -331.25

$ws.Range("H17").Value = 605.3953
$ws.Range("J17").Value = 604.0238000000001
$ws.Range("L17").Value = 1812.0714
$ws.Range("N17").Value = -2148.0714

$ws.Range("H33").Value = 77.14286
$ws.Range("I33").Value = 77.14286
$ws.Range("K33").Value = 77.14286
$ws.Range("M33").Value = 151.85714

$ws.Range("H43").Value = 2024.75
$ws.Range("I43").Value = 1466.6666
$ws.Range("J43").Value = 2359.6
$ws.Range("K43").Value = 1466.6666
$ws.Range("L43").Value = 2359.6
$ws.Range("M43").Value = -1397.6666
$ws.Range("N43").Value = -2497.6

$ws.Range("H62").Value = 2317.85
$ws.Range("I62").Value = 1601.3636
$ws.Range("J62").Value = 3193.5557
$ws.Range("K62").Value = 1601.3636
$ws.Range("L62").Value = 3193.5557
$ws.Range("M62").Value = -977.3635999999999
$ws.Range("N62").Value = -4441.5557

$ws.Range("H65").Value = 2317.85
$ws.Range("I65").Value = 1601.3636
$ws.Range("J65").Value = 3193.5557
$ws.Range("K65").Value = 8006.817999999999
$ws.Range("L65").Value = 15967.7785
$ws.Range("M65").Value = -4886.817999999999
$ws.Range("N65").Value = -22207.7785

$ws.Range("H129").Value = 1182.6621
$ws.Range("J129").Value = 1192.0137
$ws.Range("L129").Value = 3576.0411
$ws.Range("N129").Value = -13576.0411

$ws.Range("H132").Value = 4073.9
$ws.Range("I132").Value = 4905.357
$ws.Range("K132").Value = 14716.071
$ws.Range("M132").Value = -12186.071

$ws.Range("H137").Value = 79160.234
$ws.Range("I137").Value = 2295.5
$ws.Range("J137").Value = 113322.336
$ws.Range("K137").Value = 6886.5
$ws.Range("L137").Value = 339967.008
$ws.Range("M137").Value = -4336.5
$ws.Range("N137").Value = -345067.008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18942.623
$ws.Range("I32").Value = 19244.742
$ws.Range("J32").Value = 13101.667
$ws.Range("K32").Value = 19244.742
$ws.Range("L32").Value = 13101.667
$ws.Range("M32").Value = -18957.742
$ws.Range("N32").Value = -13675.667

$ws.Range("H45").Value = 3649.5715
$ws.Range("I45").Value = 3320.1875
$ws.Range("K45").Value = 3320.1875
$ws.Range("M45").Value = -2943.1875

$ws.Range("H102").Value = 1447.7858
$ws.Range("I102").Value = 1227
$ws.Range("J102").Value = 1999.75
$ws.Range("K102").Value = 1227
$ws.Range("L102").Value = 1999.75
$ws.Range("M102").Value = 395
$ws.Range("N102").Value = -5243.75

$ws.Range("H110").Value = 990
$ws.Range("I110").Value = 487.5
$ws.Range("K110").Value = 487.5
$ws.Range("M110").Value = 1557.5

$ws.Range("H132").Value = 19289.896
$ws.Range("I132").Value = 1953.9584
$ws.Range("J132").Value = 102502.4
$ws.Range("K132").Value = 5861.8752
$ws.Range("L132").Value = 307507.2
$ws.Range("M132").Value = -3331.8752
$ws.Range("N132").Value = -312567.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16569.523
$ws.Range("I31").Value = 24515.385
$ws.Range("J31").Value = 3657.5
$ws.Range("K31").Value = 24515.385
$ws.Range("L31").Value = 3657.5
$ws.Range("M31").Value = -24220.385
$ws.Range("N31").Value = -4247.5

$ws.Range("H34").Value = 16569.523
$ws.Range("I34").Value = 24515.385
$ws.Range("J34").Value = 3657.5
$ws.Range("K34").Value = 24515.385
$ws.Range("L34").Value = 3657.5
$ws.Range("M34").Value = -24313.385
$ws.Range("N34").Value = -4061.5

$ws.Range("H132").Value = 37736.668
$ws.Range("I132").Value = 44503.082
$ws.Range("J132").Value = 10671
$ws.Range("K132").Value = 133509.246
$ws.Range("L132").Value = 32013
$ws.Range("M132").Value = -130979.246
$ws.Range("N132").Value = -37073

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 146
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 162
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 486
$ws.Range("M12").Value = 23
$ws.Range("N12").Value = -832

$ws.Range("H126").Value = 3145
$ws.Range("I126").Value = 860
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 2580
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = 2360
$ws.Range("N126").Value = -39880

$ws.Range("H131").Value = 808.36365
$ws.Range("I131").Value = 633.3333
$ws.Range("J131").Value = 813.8333
$ws.Range("K131").Value = 1899.9999
$ws.Range("L131").Value = 2441.4999
$ws.Range("M131").Value = 3140.0001
$ws.Range("N131").Value = -12521.4999

$ws.Range("H132").Value = 1422.5454
$ws.Range("I132").Value = 1092.5714
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 9833.142600000001
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -7303.142600000001
$ws.Range("N132").Value = -23060

$ws.Range("H141").Value = 4841.1113
$ws.Range("I141").Value = 4765
$ws.Range("J141").Value = 4993.3335
$ws.Range("K141").Value = 14295
$ws.Range("L141").Value = 14980.0005
$ws.Range("M141").Value = -9115
$ws.Range("N141").Value = -25340.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3117.5557
$ws.Range("I97").Value = 1151.1428
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 1151.1428
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -655.1428000000001
$ws.Range("N97").Value = -10992

$ws.Range("H102").Value = 100002480
$ws.Range("I102").Value = 100002480
$ws.Range("K102").Value = 100002480
$ws.Range("M102").Value = -100000858

$ws.Range("H126").Value = 3587.158
$ws.Range("I126").Value = 2673.5386
$ws.Range("K126").Value = 8020.6158
$ws.Range("M126").Value = -5550.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4049.1738
$ws.Range("I61").Value = 1678.6666
$ws.Range("J61").Value = 6635.1816
$ws.Range("K61").Value = 1678.6666
$ws.Range("L61").Value = 6635.1816
$ws.Range("M61").Value = -1476.6666
$ws.Range("N61").Value = -7039.1816

$ws.Range("H113").Value = 4049.1738
$ws.Range("I113").Value = 1678.6666
$ws.Range("J113").Value = 6635.1816
$ws.Range("K113").Value = 1678.6666
$ws.Range("L113").Value = 6635.1816
$ws.Range("M113").Value = 491.3334
$ws.Range("N113").Value = -10975.1816

$ws.Range("H132").Value = 1836.5714
$ws.Range("I132").Value = 1271.8948
$ws.Range("K132").Value = 3815.6844
$ws.Range("M132").Value = -1285.6844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2530.1
$ws.Range("I132").Value = 1217.3334
$ws.Range("J132").Value = 4499.25
$ws.Range("K132").Value = 3652.0002
$ws.Range("L132").Value = 13497.75
$ws.Range("M132").Value = -1122.0002
$ws.Range("N132").Value = -18557.75
